# try to fix render markdown
# Replace the literal "<br>" markers inside a handful of shared strings with
# an actual line break (" " + newline) and, for the single-break headings,
# prefix the text after the break with "**test**".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$newGlobalHeading = "Global: Education, Healthcare and " + $nl + "**test**Renewable energy in low-income countries"
$newShareAllocated = "Share allocated to Global spending options " + $nl + "when 5 out of 13 options are randomly selected " + $nl + "(4 out of 13 being of Global nature)"
$newEducationHealthcare = "Global: Education and Healthcare " + $nl + "**test**in low-income countries"
$newRenewableEnergy = "Global: Renewable energy and " + $nl + "**test**infrastructure to cope with climate change"
$newLossDamage = "Global: Loss and Damage Fund (to " + $nl + "**test**rebuild after climate disasters)"

$ws.Range("E2:E12").Value = $newGlobalHeading
$ws.Range("E24:E34").Value = $newShareAllocated
$ws.Range("E35:E45").Value = $newEducationHealthcare
$ws.Range("E46:E56").Value = $newRenewableEnergy
$ws.Range("E57:E67").Value = $newLossDamage
